$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected num_matches for season 11 (row 13)
$ws.Range("E13").Value = 1182353

# Bring the bordered/centered "season index" style from A13 down onto the
# two new rows before filling them in (A14:A15)
$ws.Range("A13").Copy()
$ws.Range("A14:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 14 - season 12 (M3_01 Wolf 2021)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "M3_01 Wolf 2021"
$ws.Range("C14").Value = 9637
$ws.Range("D14").Value = 10653
$ws.Range("E14").Value = 808651
$ws.Range("F14").Value = 9916
$ws.Range("G14").Value = 10044
$ws.Range("H14").Value = 10295

# New row 15 - season 13 (M3_02 Love 2021)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "M3_02 Love 2021"
$ws.Range("C15").Value = 9684
$ws.Range("D15").Value = 10714
$ws.Range("E15").Value = 917491
$ws.Range("F15").Value = 9975
$ws.Range("G15").Value = 10097
$ws.Range("H15").Value = 10325
